$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 316, shifting existing rows 316:395 down to 318:397.
$ws.Rows("316:317").Insert()

# New row 316 data
$ws.Cells.Item(316, 1).Value2 = 4
$ws.Cells.Item(316, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(316, 3).Value2 = "Los Lagos"
$ws.Cells.Item(316, 4).Value2 = 44754
$ws.Cells.Item(316, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(316, 5).Value2 = 10
$ws.Cells.Item(316, 6).Value2 = 100114013
$ws.Cells.Item(316, 7).Value2 = "Zanahoria"
$ws.Cells.Item(316, 8).Value2 = "Sin especificar"
$ws.Cells.Item(316, 9).Value2 = "Primera"
$ws.Cells.Item(316, 10).Value2 = 450
$ws.Cells.Item(316, 11).Value2 = 13000
$ws.Cells.Item(316, 12).Value2 = 13000
$ws.Cells.Item(316, 13).Value2 = 13000
$ws.Cells.Item(316, 14).Value2 = "`$/saco 20 kilos"
$ws.Cells.Item(316, 15).Value2 = "Chillán"
$ws.Cells.Item(316, 16).Value2 = 650
$ws.Cells.Item(316, 17).Value2 = 20
$ws.Cells.Item(316, 18).Value2 = "Hortaliza"

# New row 317 data
$ws.Cells.Item(317, 1).Value2 = 4
$ws.Cells.Item(317, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(317, 3).Value2 = "Los Lagos"
$ws.Cells.Item(317, 4).Value2 = 44754
$ws.Cells.Item(317, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(317, 5).Value2 = 10
$ws.Cells.Item(317, 6).Value2 = 100114013
$ws.Cells.Item(317, 7).Value2 = "Zanahoria"
$ws.Cells.Item(317, 8).Value2 = "Sin especificar"
$ws.Cells.Item(317, 9).Value2 = "Primera"
$ws.Cells.Item(317, 10).Value2 = 450
$ws.Cells.Item(317, 11).Value2 = 10000
$ws.Cells.Item(317, 12).Value2 = 10000
$ws.Cells.Item(317, 13).Value2 = 10000
$ws.Cells.Item(317, 14).Value2 = "`$/saco 20 kilos"
$ws.Cells.Item(317, 15).Value2 = "Provincia de Llanquihue"
$ws.Cells.Item(317, 16).Value2 = 500
$ws.Cells.Item(317, 17).Value2 = 20
$ws.Cells.Item(317, 18).Value2 = "Hortaliza"
